$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): repurpose old "Recorded Split 3/4" columns (R,S)
# into four new "Recorded Finish Leg n" columns (S,T,U,V). ---
$ws.Range("R1").ClearContents()
$ws.Range("S1").Value2 = "Recorded Finish Leg 1"
$ws.Range("T1").Value2 = "Recorded Finish Leg 2"
$ws.Range("U1").Value2 = "Recorded Finish Leg 3"
$ws.Range("V1").Value2 = "Recorded Finish Leg 4"

# Old column R (rows 2-6) held the "leg-3 finish" raw reference that now
# lives in column U; clear it out completely so it is empty like the
# target layout (no leftover cell/style records).
$ws.Range("R2:R6").Clear()

# --- New "raw" recorded-finish columns S:V for rows 2-6, pulling straight
# from the detail rows below (same cells the old G/J/R/S columns used). ---
$ws.Range("S2").Formula = "=B16"
$ws.Range("T2").Formula = "=B17"
$ws.Range("U2").Formula = "=B18"
$ws.Range("V2").Formula = "=B19"

$ws.Range("S3").Formula = "=B20"
$ws.Range("T3").Formula = "=B21"
$ws.Range("U3").Formula = "=B22"
$ws.Range("V3").Formula = "=B23"

$ws.Range("S4").Formula = "=B24"
$ws.Range("T4").Formula = "=B25"
$ws.Range("U4").Formula = "=B26"
$ws.Range("V4").Formula = "=B27"

$ws.Range("S5").Formula = "=B28"
$ws.Range("T5").Formula = "=B29"
$ws.Range("U5").Formula = "=B30"
$ws.Range("V5").Formula = "=B31"

$ws.Range("S6").Formula = "=B32"
$ws.Range("T6").Formula = "=B33"
$ws.Range("U6").Formula = "=B34"
$ws.Range("V6").Formula = "=B35"

# --- Leg-1 split columns F, G now derive from the new S/T columns. ---
$ws.Range("F2").Formula = "=S2"
$ws.Range("G2").Formula = "=F2"

$ws.Range("F3:F4").Formula = "=S3"
$ws.Range("G3:G6").Formula = "=F3"

$ws.Range("F5:F6").Formula = "=S5"

# --- Leg-2 split columns I, J now derive from S/T/G. ---
$ws.Range("I2:I6").Formula = "=T2-S2"
$ws.Range("J2:J6").Formula = "=G2+I2"

# --- Leg-3 split column L now derives from T/U instead of J/R. ---
$ws.Range("L2").Formula = "=U2-MIN(L`$10,T2)"
$ws.Range("L3").Formula = "=U3-G14"
$ws.Range("L4:L6").Formula = "=U4-MIN(L`$10,T4)"

$ws.Range("M2:M6").Formula = "=J2+L2"

# --- Leg-4 split column O now derives from U/V instead of R/S. ---
$ws.Range("O2").Formula = "=V2-MIN(O`$10,U2)"
$ws.Range("O3").Formula = "=V3-MIN(O`$10,U3)"
$ws.Range("O4").Formula = "=V4-MIN(O`$10,U4)"
$ws.Range("O5").Formula = "=V5-MIN(O`$10,U5)"
$ws.Range("O6").Formula = "=V6-MIN(O`$10,U6)"

$ws.Range("P2:P3").Formula = "=M2+O2"
$ws.Range("P4").Formula = "=M4+O4"
$ws.Range("P5:P6").Formula = "=M5+O5"

# --- Selection, matching the author's final cursor position. ---
$ws.Range("R20").Select()

Write-Host "edit applied"
